# Update cryptos list with latest price/volume data from GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.645.67'
$ws.Range("E2").Value = '  -2.00%  '
$ws.Range("D3").Value = '2.629.62'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.60'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.87'
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '2.638.72'
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.05'
$ws.Range("E10").Value = '  +8.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.101'
$ws.Range("E11").Value = '  -2.06%  '
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").Value = '3.092.58'
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").Value = '58.581.62'
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.90'
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("D17").Value = '2.629.45'
$ws.Range("E17").Value = '  -0.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000133'
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.39'
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '335.31'
$ws.Range("E20").Value = '  -2.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.18'
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.55'
$ws.Range("E24").Value = '  -2.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.415'
$ws.Range("E25").Value = '  +1.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.164'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.14'
$ws.Range("E28").Value = '  -1.68%  '
$ws.Range("D29").Value = '0.0₃0740'
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.88'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.79'
$ws.Range("E33").Value = '  -0.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.52'
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.91'
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '37.13'
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("E37").Value = '  -1.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.829'
$ws.Range("E38").Value = '  -3.02%  '
$ws.Range("E39").Value = '  -3.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.815'
$ws.Range("E40").Value = '  -2.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.59'
$ws.Range("E41").Value = '  +0.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '282.01'
$ws.Range("E42").Value = '  +2.47%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.601'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("E45").Value = '  -0.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.11'
$ws.Range("E46").Value = '  +2.32%  '
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0938'
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0225'
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("D50").Value = '1.941.29'
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.46'
$ws.Range("E51").Value = '  -1.65%  '
